# Change the unit test for Sample Annotation validation to suit the new
# version of MSTemplate_Creator: the workbook used to ship three helper
# sheets (Transition_Name_Annot, ISTD_Annot, Sample_Annot) but only
# Sample_Annot is needed now, so drop the other two and leave Sample_Annot
# as the sole (and therefore active/selected) worksheet.

$wb = $excel.ActiveWorkbook

# Deleting a worksheet normally pops a confirmation dialog in real Excel;
# suppress that so the delete goes through unattended.
$excel.DisplayAlerts = $false

$wb.Worksheets("Transition_Name_Annot").Delete()
$wb.Worksheets("ISTD_Annot").Delete()

# With the other sheets gone, Sample_Annot is the only tab left - make sure
# it is the active/selected sheet.
$wb.Worksheets("Sample_Annot").Activate()

$excel.DisplayAlerts = $true
